$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$sh.Left = 164.38771823543308
$sh.Top = 355.0488188976378
$sh.Width = 391.22456692913386
$sh.Height = 88.9511834023622
